$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per the edit
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 3

$ws.Range("C3").Value = 4
$ws.Range("J3").Value = 3

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1

# Update the active selection to F8
$ws.Range("F8").Select()
